# Update Leve-profit price-tracking figures (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 116.15385
$ws.Range("I5").Value = 130.55556
$ws.Range("K5").Value = 130.55556
$ws.Range("M5").Value = -15.55556000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 435070.4
$ws.Range("I9").Value = 2500142.5
$ws.Range("K9").Value = 2500142.5
$ws.Range("M9").Value = -2499973.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3500
$ws.Range("J40").Value = 3500
$ws.Range("L40").Value = 3500
$ws.Range("N40").Value = -3850

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3221.7273
$ws.Range("I86").Value = 3159.4666
$ws.Range("J86").Value = 3355.1428
$ws.Range("K86").Value = 3159.4666
$ws.Range("L86").Value = 3355.1428
$ws.Range("M86").Value = -2036.4666
$ws.Range("N86").Value = -5601.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3221.7273
$ws.Range("I89").Value = 3159.4666
$ws.Range("J89").Value = 3355.1428
$ws.Range("K89").Value = 15797.333
$ws.Range("L89").Value = 16775.714
$ws.Range("M89").Value = -10181.333
$ws.Range("N89").Value = -28007.714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 25001.334
$ws.Range("I125").Value = 2738.25
$ws.Range("J125").Value = 42811.8
$ws.Range("K125").Value = 24644.25
$ws.Range("L125").Value = 385306.2
$ws.Range("M125").Value = -22184.25
$ws.Range("N125").Value = -390226.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3641.7585
$ws.Range("I138").Value = 1293.6428
$ws.Range("J138").Value = 5833.3335
$ws.Range("K138").Value = 3880.9284
$ws.Range("L138").Value = 17500.0005
$ws.Range("M138").Value = 1259.0716
$ws.Range("N138").Value = -27780.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1228.4
$ws.Range("I88").Value = 1037.8889
$ws.Range("J88").Value = 1335.5625
$ws.Range("K88").Value = 1037.8889
$ws.Range("L88").Value = 1335.5625
$ws.Range("M88").Value = -631.8888999999999
$ws.Range("N88").Value = -2147.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1228.4
$ws.Range("I91").Value = 1037.8889
$ws.Range("J91").Value = 1335.5625
$ws.Range("K91").Value = 1037.8889
$ws.Range("L91").Value = 1335.5625
$ws.Range("M91").Value = 366.1111000000001
$ws.Range("N91").Value = -4143.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 524.73334
$ws.Range("I110").Value = 517.0769
$ws.Range("K110").Value = 517.0769
$ws.Range("M110").Value = 1527.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1949.907
$ws.Range("I132").Value = 1993
$ws.Range("K132").Value = 5979
$ws.Range("M132").Value = -3449

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2809.5
$ws.Range("I20").Value = 3136
$ws.Range("J20").Value = 2669.5715
$ws.Range("K20").Value = 3136
$ws.Range("L20").Value = 2669.5715
$ws.Range("M20").Value = -2889
$ws.Range("N20").Value = -3163.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18491.5
$ws.Range("I82").Value = 5586.4
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 5586.4
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = -5203.4
$ws.Range("N82").Value = -40766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 18491.5
$ws.Range("I85").Value = 5586.4
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 5586.4
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = -4260.4
$ws.Range("N85").Value = -42652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37073524
$ws.Range("I86").Value = 1835.3334
$ws.Range("J86").Value = 64877290
$ws.Range("K86").Value = 1835.3334
$ws.Range("L86").Value = 64877290
$ws.Range("M86").Value = -712.3334
$ws.Range("N86").Value = -64879536

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 37073524
$ws.Range("I89").Value = 1835.3334
$ws.Range("J89").Value = 64877290
$ws.Range("K89").Value = 9176.666999999999
$ws.Range("L89").Value = 324386450
$ws.Range("M89").Value = -3560.666999999999
$ws.Range("N89").Value = -324397682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 239998.5
$ws.Range("J130").Value = 239998.5
$ws.Range("L130").Value = 239998.5
$ws.Range("N130").Value = -250038.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4935.675
$ws.Range("I31").Value = 2084.077
$ws.Range("J31").Value = 6308.6665
$ws.Range("K31").Value = 2084.077
$ws.Range("L31").Value = 6308.6665
$ws.Range("M31").Value = -1789.077
$ws.Range("N31").Value = -6898.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4935.675
$ws.Range("I34").Value = 2084.077
$ws.Range("J34").Value = 6308.6665
$ws.Range("K34").Value = 2084.077
$ws.Range("L34").Value = 6308.6665
$ws.Range("M34").Value = -1882.077
$ws.Range("N34").Value = -6712.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3040.75
$ws.Range("I99").Value = 3082.5
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 3082.5
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = -1584.5
$ws.Range("N99").Value = -5995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3040.75
$ws.Range("I126").Value = 3082.5
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 9247.5
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -6777.5
$ws.Range("N126").Value = -13937

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 195.44444
$ws.Range("J15").Value = 223.57143
$ws.Range("L15").Value = 670.71429
$ws.Range("N15").Value = -950.71429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1059.4615
$ws.Range("I34").Value = 112.4
$ws.Range("J34").Value = 1651.375
$ws.Range("K34").Value = 337.2
$ws.Range("L34").Value = 4954.125
$ws.Range("M34").Value = -253.2
$ws.Range("N34").Value = -5122.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 100
$ws.Range("I42").Value = 100
$ws.Range("K42").Value = 300
$ws.Range("M42").Value = 234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1440.1428
$ws.Range("J117").Value = 1248.8334
$ws.Range("L117").Value = 3746.5002
$ws.Range("N117").Value = -10630.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1075.1111
$ws.Range("I140").Value = 964.82355
$ws.Range("K140").Value = 2894.47065
$ws.Range("M140").Value = 2285.52935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4371.3794
$ws.Range("I80").Value = 2567
$ws.Range("K80").Value = 2567
$ws.Range("M80").Value = -1569

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4371.3794
$ws.Range("I83").Value = 2567
$ws.Range("K83").Value = 12835
$ws.Range("M83").Value = -7843

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3753.625
$ws.Range("I7").Value = 3171.6667
$ws.Range("J7").Value = 5499.5
$ws.Range("K7").Value = 3171.6667
$ws.Range("L7").Value = 5499.5
$ws.Range("M7").Value = -3059.6667
$ws.Range("N7").Value = -5723.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2660
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 2968.5715
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 2968.5715
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -3344.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3753.625
$ws.Range("I126").Value = 3171.6667
$ws.Range("J126").Value = 5499.5
$ws.Range("K126").Value = 9515.000100000001
$ws.Range("L126").Value = 16498.5
$ws.Range("M126").Value = -7045.000100000001
$ws.Range("N126").Value = -21438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 99399.60000000001
$ws.Range("J112").Value = 99399.60000000001
$ws.Range("L112").Value = 99399.60000000001
$ws.Range("N112").Value = -102353.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 913.29034
$ws.Range("I113").Value = 652.45
$ws.Range("K113").Value = 1957.35
$ws.Range("M113").Value = 212.6499999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 195000
$ws.Range("J121").Value = 195000
$ws.Range("L121").Value = 195000
$ws.Range("N121").Value = -198494
Write-Output "Goblin_Profits sheets updated"
